{"js": "// Remove the paragraph \"FABIO NUOVO ELEMENTO NEL DIZIONARIO (US)\" from the\n// bulleted list (it was deleted entirely, along with its own paragraph mark,\n// right after \"FABIO HO BISOGNO DI DEF SHARED\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetText = \"FABIO NUOVO ELEMENTO NEL DIZIONARIO (US)\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === targetText) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the paragraph \"FABIO NUOVO ELEMENTO NEL DIZIONARIO (US)\" from the\n# bulleted list (it was deleted entirely, along with its own paragraph mark,\n# right after \"FABIO HO BISOGNO DI DEF SHARED\").\n\n$d = $word.ActiveDocument\n\n$targetText = \"FABIO NUOVO ELEMENTO NEL DIZIONARIO (US)\"\n\n# Collect matching paragraphs first (collection is live; deleting while\n# forward-iterating could otherwise skip items).\n$toDelete = @()\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text.Trim() -eq $targetText) {\n        $toDelete += $p\n    }\n}\n\nforeach ($p in $toDelete) {\n    $p.Range.Delete()\n}\n"}
